$p = $ppt.ActivePresentation

# The "Conclusions" slide sits at position 14. Duplicate it, then move the
# duplicate so it sits BEFORE the original (the new slide becomes position 14
# and the original "Conclusions" slide is pushed down to position 15).
$original = $p.Slides.Item(14)
$dupRange = $original.Duplicate()
$newSlide = $dupRange.Item(1)
$newSlide.MoveTo(14)

# Fill in the content placeholder of the newly inserted slide with the note
# about the absorber-mass results (keep the "Conclusions" title as-is).
$tr = $newSlide.Shapes.Item(2).TextFrame.TextRange
$tr.Text = "Note for absorber of mass = 0.4 this is becoming quite close to the equivalent mass of building (3.3kg) so the results produced may not be reliable as mass of absorber will have "
$tr.LanguageID = "en-GB"

$r = $tr.InsertAfter("substantial ")
$r.LanguageID = "en-GB"

$r = $tr.InsertAfter("affect ")
$r.LanguageID = "en-GB"

$r = $tr.InsertAfter("on ")
$r.LanguageID = "en-GB"

$r = $tr.InsertAfter("building")
$r.LanguageID = "en-GB"
